$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns for data rows are treated as text,
# matching the source data which stores these as inline strings (not numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '34.601.48'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '1.816.80'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '228.38'
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").Value = '0.559'
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = '34.83'
$ws.Range("E8").Value = '  +7.68%  '
$ws.Range("D9").Value = '0.301'
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("D10").Value = '0.0696'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").Value = '0.0953'
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '2.076.44'
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").Value = '11.34'
$ws.Range("E13").Value = '  +2.64%  '
$ws.Range("D14").Value = '1.810.51'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").Value = '0.644'
$ws.Range("E15").Value = '  +2.73%  '
$ws.Range("D16").Value = '34.628.45'
$ws.Range("E16").Value = '  +1.04%  '
$ws.Range("D17").Value = '4.33'
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("D18").Value = '69.15'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").Value = '247.53'
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").Value = '0.0₃0802'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '11.56'
$ws.Range("E21").Value = '  +5.58%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '4.21'
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").Value = '172.21'
$ws.Range("E24").Value = '  +6.07%  '
$ws.Range("D25").Value = '2.09'
$ws.Range("E25").Value = '  +1.96%  '
$ws.Range("D26").Value = '7.47'
$ws.Range("E26").Value = '  +3.98%  '
$ws.Range("D27").Value = '16.76'
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("D28").Value = '0.117'
$ws.Range("E28").Value = '  +1.28%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '4.03'
$ws.Range("E30").Value = '  +4.73%  '
$ws.Range("D31").Value = '0.0533'
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("D32").Value = '3.86'
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("D34").Value = '1.86'
$ws.Range("E34").Value = '  +2.71%  '
$ws.Range("D35").Value = '2.61'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '1.421.01'
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("D37").Value = '0.677'
$ws.Range("E37").Value = '  +2.40%  '
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").Value = '0.0192'
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").Value = '86.13'
$ws.Range("E40").Value = '  +4.51%  '
$ws.Range("D41").Value = '2.86'
$ws.Range("E41").Value = '  +4.40%  '
$ws.Range("D42").Value = '0.964'
$ws.Range("E42").Value = '  +4.39%  '
$ws.Range("D43").Value = '2.41'
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("D44").Value = '13.92'
$ws.Range("E44").Value = '  -1.57%  '
$ws.Range("D45").Value = '0.0525'
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("D46").Value = '1.11'
$ws.Range("E46").Value = '  +3.08%  '
$ws.Range("D47").Value = '6.12'
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("D48").Value = '1.978.84'
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("D49").Value = '106.17'
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("D50").Value = '0.0₆0131'
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.18%  '
